$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: update the existing "manager" account in place (row 1 at this point)
$ws.Range("A1").Value = "manager"
$ws.Range("B1").Value = "48xaI2RykI2DkQrmO0hchQ=="
$ws.Range("C1").Value = "Manajer Utama"

# Step 2: append the additional accounts below it
$ws.Range("A2").Value = "rizki"
$ws.Range("B2").Value = "4F8Q91ePWMvdNEbFmKiQqA=="
$ws.Range("C2").Value = "Manajer Utama"

$ws.Range("A3").Value = "dan"
$ws.Range("B3").Value = "hOUNFYHr3v23KAwmJanNMQ=="
$ws.Range("C3").Value = "Manajer Utama"

# Step 3: insert a header row above the data
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "Nama"
$ws.Range("B1").Value = "password"
$ws.Range("C1").Value = "jabatan"

# Step 4: append the last account row
$ws.Range("A5").Value = "riz"
$ws.Range("B5").Value = "Gby8QXqKlgnMNWi0Z5xgbw=="
$ws.Range("C5").Value = "Manajer Utama"

# Column sizing for the password / jabatan columns
$ws.Columns.Item(2).ColumnWidth = 27.5
$ws.Columns.Item(3).ColumnWidth = 13

# Match the resulting selection (column B selected)
$ws.Range("B1:B1048576").Select()
